# Update automàtic: dades i banners [2026-02-05 07:39]
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New extraction timestamps (DATA_EXTRACCIO, column E) for rows 2-36
$timestamps = @{
    2  = "2026-02-05 07:38:13"
    3  = "2026-02-05 07:38:15"
    4  = "2026-02-05 07:38:17"
    5  = "2026-02-05 07:38:20"
    6  = "2026-02-05 07:38:22"
    7  = "2026-02-05 07:38:25"
    8  = "2026-02-05 07:38:27"
    9  = "2026-02-05 07:38:30"
    10 = "2026-02-05 07:38:32"
    11 = "2026-02-05 07:38:34"
    12 = "2026-02-05 07:38:37"
    13 = "2026-02-05 07:38:39"
    14 = "2026-02-05 07:38:42"
    15 = "2026-02-05 07:38:44"
    16 = "2026-02-05 07:38:46"
    17 = "2026-02-05 07:38:49"
    18 = "2026-02-05 07:38:52"
    19 = "2026-02-05 07:38:54"
    20 = "2026-02-05 07:38:57"
    21 = "2026-02-05 07:38:59"
    22 = "2026-02-05 07:39:01"
    23 = "2026-02-05 07:39:03"
    24 = "2026-02-05 07:39:06"
    25 = "2026-02-05 07:39:08"
    26 = "2026-02-05 07:39:11"
    27 = "2026-02-05 07:39:13"
    28 = "2026-02-05 07:39:16"
    29 = "2026-02-05 07:39:18"
    30 = "2026-02-05 07:39:21"
    31 = "2026-02-05 07:39:23"
    32 = "2026-02-05 07:39:25"
    33 = "2026-02-05 07:39:28"
    34 = "2026-02-05 07:39:30"
    35 = "2026-02-05 07:39:33"
    36 = "2026-02-05 07:39:35"
}

foreach ($row in $timestamps.Keys) {
    $ws.Cells.Item($row, 5).Value = $timestamps[$row]
}

# Row 30 (ZC - Setcases - Ulldeter): banner data now populated
$ws.Cells.Item(30, 7).Value  = "56 cm"
$ws.Cells.Item(30, 8).Value  = "'56%"
$ws.Cells.Item(30, 9).Value  = "0.0 mm"
$ws.Cells.Item(30, 11).Value = "0.0 MJ/m2"
$ws.Cells.Item(30, 12).Value = "15.8 km/h - 291º 0:52 TU"
$ws.Cells.Item(30, 13).Value = "-1.6 °C 1:10 TU"
$ws.Cells.Item(30, 14).Value = "-3.9 °C 0:10 TU"
$ws.Cells.Item(30, 15).Value = "-2.7 °C"

# Row 35 (VS - Vielha e Mijaran - Lac Redon): banner data now populated
$ws.Cells.Item(35, 7).Value  = "199 cm"
$ws.Cells.Item(35, 8).Value  = "'94%"
$ws.Cells.Item(35, 9).Value  = "0.0 mm"
$ws.Cells.Item(35, 11).Value = "0.0 MJ/m2"
$ws.Cells.Item(35, 12).Value = "0.0 km/h - 0º 0:00 TU"
$ws.Cells.Item(35, 13).Value = "-3.4 °C 0:16 TU"
$ws.Cells.Item(35, 14).Value = "-4.5 °C 3:11 TU"
$ws.Cells.Item(35, 15).Value = "-3.9 °C"
